$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Routes")

$ws.Range("A21").Value = "EuropeanWings"
$ws.Range("B21").Value = "Hungary-Budapest-Listz"
$ws.Range("C21").Value = "LHBP"
$ws.Range("D21").Value = "Paris-Beauvais-Tille"
$ws.Range("E21").Value = "LFOB"

$ws.Range("B21").VerticalAlignment = -4108

$ws.Range("D21:E21").Select()
